$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The blank separator row (row 2) was removed, shifting the SABANA..INVERTEBRATE
# block up by one row.
$ws.Rows(2).Delete()

# The PORTADA and "NICE TO MEET YOU" entries (which, after the shift above, land
# on rows 13 and 14) were removed from the sheet, leaving those rows blank while
# keeping their existing cell formatting.
$ws.Range("A13:C14").ClearContents()

# Reflect where the author's selection ended up after the edit.
$ws.Range("B20").Select()
